# Sprint 17 (changes suggested by Ronak in mail)
# - Add 5 new transaction rows (10-14) surrounded by a thin border
# - Add a TOTAL row (15) that sums the new "Credit" values
# - Update the "Ledger Generation Date" (B5) value

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Apply a thin border (all sides) to the new data block A10:E14.
#    Doing this before we touch any other cell style makes sure the new
#    border/style entries land first in the style table.
# ---------------------------------------------------------------------
$dataRange = $ws.Range("A10:E14")
$dataRange.Borders.LineStyle = 1   # xlContinuous
$dataRange.Borders.Weight = 2      # xlThin

# Keep column A formatted as text so the date-like strings we are about
# to enter are not auto-converted into date serial numbers.
$colA = $ws.Range("A10:A14")
$colA.NumberFormat = "@"

# ---------------------------------------------------------------------
# 2) Populate the new rows
# ---------------------------------------------------------------------
$rows = @(
    @("2020-10-03", "TEST", "IN-0005", "INVOICE", 0),
    @("2020-10-03", "TEST", "IN-0006", "INVOICE", 0),
    @("2020-10-03", "TEST", "IN-0006", "INVOICE", 0),
    @("2020-10-03", "TEST", "IN-0006", "INVOICE", 0),
    @("2020-10-03", "TEST", "IN-0009", "INVOICE", 5)
)

$r = 10
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

# Re-assert the plain thin-border style on column A (writing the text
# values above can perturb the style chosen for date-like text), and
# restore the thin border around the whole block just in case.
$colA.Style = "Normal"
$dataRange.Borders.LineStyle = 1
$dataRange.Borders.Weight = 2

# ---------------------------------------------------------------------
# 3) Add the TOTAL row (row 15). Set the value/formula BEFORE touching
#    the formatting of these cells, otherwise the formula engine can
#    end up with a stale cached result for the SUM formula.
# ---------------------------------------------------------------------
$ws.Range("D15").Value = "TOTAL"
$ws.Range("E15").Formula = "=SUM(E10:E14)"

# Reuse the header style (bold white text on dark fill with a thick
# border) from D9/E9 for the TOTAL row.
$ws.Range("D9:E9").Copy() | Out-Null
$ws.Range("D15:E15").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Make sure the value/formula are still correct after the format paste.
$ws.Range("D15").Value = "TOTAL"
$ws.Range("E15").Formula = "=SUM(E10:E14)"

# ---------------------------------------------------------------------
# 4) Update the Ledger Generation Date in B5
# ---------------------------------------------------------------------
$b5 = $ws.Range("B5")
$b5.NumberFormat = "@"
$b5.Value = "2020-10-05"
$b5.Style = "Normal"

# ---------------------------------------------------------------------
# 5) Recalculate so the SUM formula has a fresh cached value.
# ---------------------------------------------------------------------
$excel.CalculateFullRebuild()
